# Update positional_accuracy and other details (IGC_sample data)
# Applies the edited values for rows 2-37 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  A = 1063; C = 50;  D = "bella";        E = "bella"; F = 1; G = 1 },
    @{ Row = 3;  A = 545;  C = 29;  D = "veraz";         E = "veraz"; F = 1; G = 1 },
    @{ Row = 4;  A = 1107; C = 94;  D = "facha";         E = "facha"; F = 1; G = 1 },
    @{ Row = 5;  A = 573;  C = 57;  D = "leña";          E = "leña"; F = 1; G = 1 },
    @{ Row = 6;  A = 627;  C = 111; D = "raya";          E = "ya, ca, raya"; F = 0; G = 1 },
    @{ Row = 7;  A = 558;  C = 42;  D = "patria";        E = "patria"; F = 1; G = 1 },
    @{ Row = 8;  A = 527;  C = 11;  D = "labor";         E = "labor"; F = 1; G = 1 },
    @{ Row = 9;  A = 522;  C = 6;   D = "veloz";         E = "ver, lo, feloz, lelo, veloz"; F = 0; G = 1 },
    @{ Row = 10; A = 1084; C = 71;  D = "pina";          E = "pilas, pilas, pinas"; F = 0; G = 1 },
    @{ Row = 11; A = 668;  C = 32;  D = "vigente";       E = "vigente"; F = 1; G = 1 },
    @{ Row = 12; A = 672;  C = 36;  D = "caída";         E = "caída"; F = 1; G = 1 },
    @{ Row = 13; A = 1226; C = 93;  D = "carajo";        E = "garajo, carajo"; F = 0; G = 1 },
    @{ Row = 14; A = 647;  C = 11;  D = "aspecto";       E = "aspecto"; F = 1; G = 1 },
    @{ Row = 15; A = 690;  C = 54;  D = "rechazo";       E = "rechazo"; F = 1; G = 1 },
    @{ Row = 16; A = 728;  C = 92;  D = "concorde";      E = "concorde"; F = 1; G = 1 },
    @{ Row = 17; A = 1253; C = 120; D = "judía";         E = "judía"; F = 1; G = 1 },
    @{ Row = 18; A = 710;  C = 74;  D = "almuerzo";      E = "almuerzo"; F = 1; G = 1 },
    @{ Row = 19; A = 756;  C = 120; D = "judía";         E = "judía"; F = 1; G = 1 },
    @{ Row = 20; A = 940;  C = 144; D = "distinto";      E = "distinto"; F = 1; G = 1 },
    @{ Row = 21; A = 1332; C = 79;  D = "preparatoria";  E = "pre, pe, prepar, preparotoria, preparatoria"; F = 0; G = 1 },
    @{ Row = 22; A = 840;  C = 44;  D = "inundación";    E = "unun, unincia, unin, unenzación, undiz, undi, uni, unde, indi, ninunación, indunación, indunaz, indininun, indi, indu, indunización"; F = 0; G = 0 },
    @{ Row = 23; A = 1390; C = 137; D = "carnicería";    E = "carnifería"; F = 0; G = 0 },
    @{ Row = 24; A = 1420; C = 167; D = "volcán";        E = "volcán"; F = 1; G = 1 },
    @{ Row = 25; A = 1289; C = 36;  D = "economía";      E = "economía"; F = 1; G = 1 },
    @{ Row = 26; A = 815;  C = 19;  D = "esperma";       E = "esperma"; F = 1; G = 1 },
    @{ Row = 27; A = 924;  C = 128; D = "emisión";       E = "emisión"; F = 1; G = 1 },
    @{ Row = 28; A = 994;  C = 198; D = "superioridad";  E = "superioridad"; F = 1; G = 1 },
    @{ Row = 29; A = 1835; C = 54;  D = "quejido";       E = "tejido"; F = 0; G = 0 },
    @{ Row = 30; A = 1856; C = 24;  D = "destello";      E = "destello"; F = 1; G = 1 },
    @{ Row = 31; A = 1862; C = 67;  D = "tocador";       E = "tocador"; F = 1; G = 1 },
    @{ Row = 32; A = 1865; C = 39;  D = "léxico";        E = "léstico"; F = 0; G = 0 },
    @{ Row = 33; A = 1801; C = 44;  D = "merluza";       E = "merluza"; F = 1; G = 1 },
    @{ Row = 34; A = 1836; C = 71;  D = "ventanal";      E = "ventanal"; F = 1; G = 1 },
    @{ Row = 35; A = 1852; C = 51;  D = "pilares";       E = "pilares"; F = 1; G = 1 },
    @{ Row = 36; A = 1849; C = 21;  D = "delirio";       E = "delirio"; F = 1; G = 1 },
    @{ Row = 37; A = 1859; C = 6;   D = "capellán";      E = "capellán"; F = 1; G = 1 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = $r.A
    $ws.Cells.Item($i, 3).Value = $r.C
    $ws.Cells.Item($i, 4).Value = $r.D
    $ws.Cells.Item($i, 5).Value = $r.E
    $ws.Cells.Item($i, 6).Value = $r.F
    $ws.Cells.Item($i, 7).Value = $r.G
}
